# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 3: Intel(R) Wi-Fi 6 AX200 160MHz - 23.110.0.5
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 58
$ws.Range("D3").Value = 98.09999999999999

# Row 4: Totals
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 58

# Row 14: Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8
$ws.Range("B14").Value = 331283
